# Applies the "personal_file" update:
#  - Row 4: Lfmc_id/Passport_series/Passport_id/Phone_number become numeric;
#           Birthday_place changes from "Москва" to "Екатерингбург"
#  - Row 7: values replaced, and become text (inlineStr) instead of numeric
#  - Row 8: values replaced; Lfmc_id/Passport_series/Passport_id become numeric
#  - Rows 9 and 10: brand-new family-member records appended

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($ws, $addr, $val) {
    # These cells are already General-formatted, so a plain numeric
    # assignment is stored as a real number without touching styles.
    $c = $ws.Range($addr)
    $c.Value = $val
}

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    # Force the cell to stay text even though the value looks numeric,
    # matching Excel's own "Text" input behaviour, then drop the format
    # override again so no extra style survives on the cell.
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-NumCell $ws "A4" 3
Set-NumCell $ws "B4" 1111
Set-NumCell $ws "C4" 285463
Set-TextCell $ws "D4" "Екатерингбург"
Set-TextCell $ws "E4" "Дефолт Сити, улица Пушкина, дом 14"
Set-TextCell $ws "F4" "Холост"
Set-TextCell $ws "G4" "САФУ"
Set-TextCell $ws "H4" "Работа пример"
Set-TextCell $ws "I4" "-"
Set-NumCell $ws "J4" 88005553537
Set-TextCell $ws "K4" "Годен, категория А. Подлежит призыву к службе в ВС РФ"
Set-TextCell $ws "A7" "6"
Set-TextCell $ws "B7" "1234"
Set-TextCell $ws "C7" "123456"
Set-TextCell $ws "D7" "ahgdsfhgsd"
Set-TextCell $ws "E7" "aaaaa"
Set-TextCell $ws "F7" "ahgdsfhwdfiuhvdngsd"
Set-TextCell $ws "G7" "nfwebcxsow"
Set-TextCell $ws "H7" "uwdhcvbsnkas"
Set-TextCell $ws "I7" "mcxqucnxjanksjweidjsbancjsdn"
Set-TextCell $ws "J7" "saduqwbdas"
Set-TextCell $ws "K7" "mcxqucnxjanksjweidjsbancjsdn"
Set-NumCell $ws "A8" 10
Set-NumCell $ws "B8" 1111
Set-NumCell $ws "C8" 123123
Set-TextCell $ws "D8" "gwerewrfsfg"
Set-TextCell $ws "E8" "whfdgsdf"
Set-TextCell $ws "F8" "vcxht4uf"
Set-TextCell $ws "G8" "vcxfasfs"
Set-TextCell $ws "H8" "sdgdfijcjmvs"
Set-TextCell $ws "I8" "dgduhsdcjsnskdf"
Set-TextCell $ws "J8" "vjsidwhefsjs"
Set-TextCell $ws "K8" "fsudhfsdjnasjsni"
Set-NumCell $ws "A9" 11
Set-NumCell $ws "B9" 1212
Set-NumCell $ws "C9" 151234
Set-TextCell $ws "D9" "fsdgsdhsj"
Set-TextCell $ws "E9" "sfkghljfkyurj"
Set-TextCell $ws "F9" "fmndgfyst"
Set-TextCell $ws "G9" "gdkeyjhdfh"
Set-TextCell $ws "H9" "jetykgdhgfj"
Set-TextCell $ws "I9" "jstrbvssgh"
Set-TextCell $ws "J9" "dtyijdseg"
Set-TextCell $ws "K9" "gfhnfgjrstt"
Set-NumCell $ws "A10" 9
Set-NumCell $ws "B10" 2323
Set-NumCell $ws "C10" 124167
Set-TextCell $ws "D10" "gadfgb"
Set-TextCell $ws "E10" "Moscow"
Set-TextCell $ws "F10" "bxcvbsdfvbx"
Set-TextCell $ws "G10" "xvcbsdfg"
Set-TextCell $ws "H10" "bvxdfagb"
Set-TextCell $ws "I10" "asfhhhagchawbchanwjenfj"
Set-TextCell $ws "J10" "xdfsgg"
Set-TextCell $ws "K10" "asfhhhagchawbchanwjenfj"
